$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H12").Value = 104.90909
$ws_ALC.Range("I12").Value = 113.666664
$ws_ALC.Range("J12").Value = 65.5
$ws_ALC.Range("K12").Value = 113.666664
$ws_ALC.Range("L12").Value = 65.5
$ws_ALC.Range("M12").Value = 56.333336
$ws_ALC.Range("N12").Value = -405.5

$ws_ALC.Range("H62").Value = 7333
$ws_ALC.Range("I62").Value = 6666
$ws_ALC.Range("K62").Value = 6666
$ws_ALC.Range("M62").Value = -6042

$ws_ALC.Range("H65").Value = 7333
$ws_ALC.Range("I65").Value = 6666
$ws_ALC.Range("K65").Value = 33330
$ws_ALC.Range("M65").Value = -30210

$ws_ALC.Range("H98").Value = 982
$ws_ALC.Range("I98").Value = 1102.5
$ws_ALC.Range("J98").Value = 500
$ws_ALC.Range("K98").Value = 1102.5
$ws_ALC.Range("L98").Value = 500
$ws_ALC.Range("M98").Value = 395.5
$ws_ALC.Range("N98").Value = -3496

$ws_ALC.Range("H100").Value = 1837.091
$ws_ALC.Range("I100").Value = 1420
$ws_ALC.Range("K100").Value = 1420
$ws_ALC.Range("M100").Value = -879

$ws_ALC.Range("H122").Value = 982
$ws_ALC.Range("I122").Value = 1102.5
$ws_ALC.Range("J122").Value = 500
$ws_ALC.Range("K122").Value = 3307.5
$ws_ALC.Range("L122").Value = 1500
$ws_ALC.Range("M122").Value = -857.5
$ws_ALC.Range("N122").Value = -6400

$ws_ALC.Range("H129").Value = 5209.5293
$ws_ALC.Range("J129").Value = 6542.775
$ws_ALC.Range("L129").Value = 19628.325
$ws_ALC.Range("N129").Value = -29628.325

$ws_ALC.Range("H135").Value = 1171.5
$ws_ALC.Range("I135").Value = 816.92
$ws_ALC.Range("J135").Value = 10036
$ws_ALC.Range("K135").Value = 7352.28
$ws_ALC.Range("L135").Value = 90324
$ws_ALC.Range("M135").Value = -4817.28
$ws_ALC.Range("N135").Value = -95394

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H45").Value = 1737.75
$ws_ARM.Range("I45").Value = 1804.1538
$ws_ARM.Range("K45").Value = 1804.1538
$ws_ARM.Range("M45").Value = -1427.1538

$ws_ARM.Range("H61").Value = 3004.625
$ws_ARM.Range("I61").Value = 2774.4546
$ws_ARM.Range("K61").Value = 2774.4546
$ws_ARM.Range("M61").Value = -2562.4546

$ws_ARM.Range("H110").Value = 2013.2
$ws_ARM.Range("I110").Value = 1647.4286
$ws_ARM.Range("J110").Value = 2866.6667
$ws_ARM.Range("K110").Value = 1647.4286
$ws_ARM.Range("L110").Value = 2866.6667
$ws_ARM.Range("M110").Value = 397.5714
$ws_ARM.Range("N110").Value = -6956.6667

$ws_ARM.Range("H136").Value = 3004.625
$ws_ARM.Range("I136").Value = 2774.4546
$ws_ARM.Range("K136").Value = 8323.363799999999
$ws_ARM.Range("M136").Value = -5773.363799999999

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H86").Value = 5407708
$ws_BSM.Range("I86").Value = 9093099
$ws_BSM.Range("K86").Value = 9093099
$ws_BSM.Range("M86").Value = -9091976

$ws_BSM.Range("H89").Value = 5407708
$ws_BSM.Range("I89").Value = 9093099
$ws_BSM.Range("K89").Value = 45465495
$ws_BSM.Range("M89").Value = -45459879

$ws_BSM.Range("H99").Value = 2287.8462
$ws_BSM.Range("I99").Value = 2051.6667
$ws_BSM.Range("J99").Value = 2490.2856
$ws_BSM.Range("K99").Value = 2051.6667
$ws_BSM.Range("L99").Value = 2490.2856
$ws_BSM.Range("M99").Value = -553.6667000000002
$ws_BSM.Range("N99").Value = -5486.2856

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 797.617
$ws_CUL.Range("I5").Value = 480.77274
$ws_CUL.Range("J5").Value = 1076.44
$ws_CUL.Range("K5").Value = 1442.31822
$ws_CUL.Range("L5").Value = 3229.32
$ws_CUL.Range("M5").Value = -1330.31822
$ws_CUL.Range("N5").Value = -3453.32

$ws_CUL.Range("H68").Value = 1497.5
$ws_CUL.Range("J68").Value = 1863.3334
$ws_CUL.Range("L68").Value = 5590.0002
$ws_CUL.Range("N68").Value = -7212.0002

$ws_CUL.Range("H71").Value = 1497.5
$ws_CUL.Range("J71").Value = 1863.3334
$ws_CUL.Range("L71").Value = 16770.0006
$ws_CUL.Range("N71").Value = -24882.0006

$ws_CUL.Range("H105").Value = 24603.143
$ws_CUL.Range("J105").Value = 24603.143
$ws_CUL.Range("L105").Value = 73809.429
$ws_CUL.Range("N105").Value = -79051.429

$ws_CUL.Range("H131").Value = 3509699.8
$ws_CUL.Range("I131").Value = 66666664
$ws_CUL.Range("J131").Value = 979.44446
$ws_CUL.Range("K131").Value = 199999992
$ws_CUL.Range("L131").Value = 2938.33338
$ws_CUL.Range("M131").Value = -199994952
$ws_CUL.Range("N131").Value = -13018.33338

$ws_CUL.Range("H135").Value = 797.617
$ws_CUL.Range("I135").Value = 480.77274
$ws_CUL.Range("J135").Value = 1076.44
$ws_CUL.Range("K135").Value = 4326.95466
$ws_CUL.Range("L135").Value = 9687.960000000001
$ws_CUL.Range("M135").Value = -1791.95466
$ws_CUL.Range("N135").Value = -14757.96

$ws_CUL.Range("H137").Value = 3478.2068
$ws_CUL.Range("J137").Value = 5386.8125
$ws_CUL.Range("L137").Value = 16160.4375
$ws_CUL.Range("N137").Value = -26360.4375

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H122").Value = 1899.2858
$ws_GSM.Range("I122").Value = 1759
$ws_GSM.Range("K122").Value = 5277
$ws_GSM.Range("M122").Value = -2827

$ws_GSM.Range("H132").Value = 3025.1538
$ws_GSM.Range("I132").Value = 2335.75
$ws_GSM.Range("K132").Value = 7007.25
$ws_GSM.Range("M132").Value = -4477.25

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 1938.1666
$ws_LTW.Range("I7").Value = 1541.5714
$ws_LTW.Range("K7").Value = 1541.5714
$ws_LTW.Range("M7").Value = -1429.5714

$ws_LTW.Range("H32").Value = 5500
$ws_LTW.Range("I32").Value = 1000
$ws_LTW.Range("J32").Value = 10000
$ws_LTW.Range("K32").Value = 1000
$ws_LTW.Range("L32").Value = 10000
$ws_LTW.Range("M32").Value = -683
$ws_LTW.Range("N32").Value = -10634

$ws_LTW.Range("H40").Value = 1885.7142
$ws_LTW.Range("I40").Value = 1360
$ws_LTW.Range("J40").Value = 3200
$ws_LTW.Range("K40").Value = 1360
$ws_LTW.Range("L40").Value = 3200
$ws_LTW.Range("M40").Value = -1224
$ws_LTW.Range("N40").Value = -3472

$ws_LTW.Range("H122").Value = 5156.3477
$ws_LTW.Range("I122").Value = 5219.8
$ws_LTW.Range("J122").Value = 4733.3335
$ws_LTW.Range("K122").Value = 15659.4
$ws_LTW.Range("L122").Value = 14200.0005
$ws_LTW.Range("M122").Value = -13209.4
$ws_LTW.Range("N122").Value = -19100.0005

$ws_LTW.Range("H126").Value = 1938.1666
$ws_LTW.Range("I126").Value = 1541.5714
$ws_LTW.Range("K126").Value = 4624.7142
$ws_LTW.Range("M126").Value = -2154.7142

$ws_LTW.Range("H132").Value = 14294021
$ws_LTW.Range("I132").Value = 22738568
$ws_LTW.Range("J132").Value = 3249.4614
$ws_LTW.Range("K132").Value = 68215704
$ws_LTW.Range("L132").Value = 9748.3842
$ws_LTW.Range("M132").Value = -68213174
$ws_LTW.Range("N132").Value = -14808.3842

$ws_LTW.Range("H136").Value = 5774.4243
$ws_LTW.Range("I136").Value = 8117.1113
$ws_LTW.Range("J136").Value = 2963.2
$ws_LTW.Range("K136").Value = 24351.3339
$ws_LTW.Range("L136").Value = 8889.599999999999
$ws_LTW.Range("M136").Value = -21801.3339
$ws_LTW.Range("N136").Value = -13989.6

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H81").Value = 71431350
$ws_WVR.Range("I81").Value = 250003340
$ws_WVR.Range("J81").Value = 2550
$ws_WVR.Range("K81").Value = 500006680
$ws_WVR.Range("L81").Value = 5100
$ws_WVR.Range("M81").Value = -500005619
$ws_WVR.Range("N81").Value = -7222

$ws_WVR.Range("H84").Value = 71431350
$ws_WVR.Range("I84").Value = 250003340
$ws_WVR.Range("J84").Value = 2550
$ws_WVR.Range("K84").Value = 2500033400
$ws_WVR.Range("L84").Value = 25500
$ws_WVR.Range("M84").Value = -2500028096
$ws_WVR.Range("N84").Value = -36108

$ws_WVR.Range("H122").Value = 1057.3
$ws_WVR.Range("I122").Value = 844.087
$ws_WVR.Range("J122").Value = 1757.8572
$ws_WVR.Range("K122").Value = 2532.261
$ws_WVR.Range("L122").Value = 5273.571599999999
$ws_WVR.Range("M122").Value = -82.26099999999997
$ws_WVR.Range("N122").Value = -10173.5716

$ws_WVR.Range("H126").Value = 1666.6666
$ws_WVR.Range("I126").Value = 1666.6666
$ws_WVR.Range("K126").Value = 4999.9998
$ws_WVR.Range("M126").Value = -2529.9998

$ws_WVR.Range("H136").Value = 5773.4165
$ws_WVR.Range("I136").Value = 878.875
$ws_WVR.Range("J136").Value = 15562.5
$ws_WVR.Range("K136").Value = 2636.625
$ws_WVR.Range("L136").Value = 46687.5
$ws_WVR.Range("M136").Value = -86.625
$ws_WVR.Range("N136").Value = -51787.5
